$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Require angle brackets for @base and @prefix values:
# D1 was "http://sales.data/purchases/2015" -> "<http://sales.data/purchases/2015>"
# D2 was "http://sales.data/purchases#"      -> "<http://sales.data/purchases#>"
# D3 was "http://sales.data/schema#"          -> "<http://sales.data/schema#>"
$ws.Range("D1").Value = "<http://sales.data/purchases/2015>"
$ws.Range("D2").Value = "<http://sales.data/purchases#>"
$ws.Range("D3").Value = "<http://sales.data/schema#>"

# Update the selection shown in the sheet view
$ws.Range("D1:D3").Select()
